# "Add Project to main branch"
# - Adds a new "national-id " column (F) with per-row id values.
# - Re-types the existing "Mobile " id column (D) from numeric ids to
#   free-text ids (several rows now carry stray characters/spaces that
#   aren't valid numbers any more - hence the text/"@" number format).
# - Fixes up a couple of stray/placeholder values that were sitting in the
#   test sheet (B3 "M" -> "Male ", A7 blank -> real text, E8 "lol " -> junk
#   test string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- column D: switch to text ids -----------------------------------
$ws.Range("D2:D8").NumberFormat = "@"
$ws.Range("D2").Value = " 0101001001"
$ws.Range("D3").Value = "0101001002 "
$ws.Range("D4").Value = "+20101001003"
$ws.Range("D5").Value = "e0101001004"
$ws.Range("D6").Value = "0101001005"
$ws.Range("D7").Value = "0101001006"
$ws.Range("D8").Value = "0101001007"

# ---- misc value fixes -------------------------------------------------
$ws.Range("B3").Value = "Male "
$ws.Range("A7").Value = "sdf sjdfb jdhfs dfsn isdfhs isf  juseif lis fsi"
$ws.Range("E8").Value = "515*/**/*-/*/*/* "

# ---- new column F: national-id ----------------------------------------
$ws.Range("F1:F8").NumberFormat = "@"
$ws.Range("F1").Value = "national-id "
$ws.Range("F2").Value = "wedwhehhwhd"
$ws.Range("F3").Value = "302022022548"
$ws.Range("F4").Value = "30256655485265"
$ws.Range("F5").Value = "30256655485265"
$ws.Range("F6").Value = " 30256655485265"
$ws.Range("F7").Value = "c                                "
$ws.Range("F8").Value = "30256455485265"

# Header cell gets a small thin box (left/right) border + smaller font,
# like the rest of the header row but distinguished as the newly added one.
$ws.Range("F1").Borders.Item(7).LineStyle = 1
$ws.Range("F1").Borders.Item(10).LineStyle = 1
$ws.Range("F1").Font.Size = 8

# Column A widened to fit the long free-text entries now living in it.
$ws.Columns.Item(1).ColumnWidth = 58.3

# Leave the cursor on the newly added header cell, like the source file.
$ws.Range("F1").Select() | Out-Null
